$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H (shifts the old "center std*" column H -> I,
# and the old "time" column I -> J), matching the diff's column layout:
#   G = center Err*   (unchanged)
#   H = center err mm (new)
#   I = center std*   (was H)
#   J = time          (was I)
#   K = time per img  (new, appended)
[void]$ws.Columns("H").Insert()

# Write the new header labels in the same order the original workbook's
# shared-string table lists them ("time per img" then "center err mm"),
# even though "time per img" ends up living in the later column K.
$ws.Range("K1").Value = "time per img"
$ws.Range("H1").Value = "center err mm"

# New "center err mm" column: per-row formula + summary stats (converts
# the pixel-space center std in column G to millimeters).
$ws.Range("H2:H7").Formula = "=G2*1670/((914.885+917.224)/2)"
$ws.Range("H8").Formula = "=AVERAGE(H2:H7)"
$ws.Range("H9").Formula = "=STDEV.S(H2:H7)"

# New "time per img" column appended at the end: per-row formula + summary
# stats (converts the total time in column J, across 99 images, into
# milliseconds per image).
$ws.Range("K2:K7").Formula = "=J2/99*1000"
$ws.Range("K8").Formula = "=AVERAGE(K2:K7)"
$ws.Range("K9").Formula = "=STDEV.S(K2:K7)"

# Match the saved selection state from the diff.
[void]$ws.Range("K8").Select()
